# Adding Answers till Q16
# ------------------------------------------------------------------
# This script is run against the document via the iron_native COM
# shim. $word / $doc / $app resolve to the host; the document itself
# is available as $word.ActiveDocument.
# ------------------------------------------------------------------

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Bold + bold-complex-script + dark red (C00000), matching the rest of
# the answers already present in the document.
$answerRprOpen  = '<w:rPr><w:b/><w:bCs/><w:color w:val="C00000"/></w:rPr>'

function New-AnswerParagraphXml([string[]]$runTexts) {
    # Builds a <w:p> whose paragraph mark AND every run carry the bold /
    # red "answer" formatting, splitting the visible text across one
    # <w:r> per entry in $runTexts (mirrors how Word splits runs when a
    # sentence is typed/edited in pieces).
    $xml = '<w:p ' + $wNs + '>'
    $xml += '<w:pPr>' + $answerRprOpen + '</w:pPr>'
    foreach ($t in $runTexts) {
        $xml += '<w:r>' + $answerRprOpen
        if ($t.StartsWith(" ") -or $t.EndsWith(" ")) {
            $xml += '<w:t xml:space="preserve">' + $t + '</w:t>'
        } else {
            $xml += '<w:t>' + $t + '</w:t>'
        }
        $xml += '</w:r>'
    }
    $xml += '</w:p>'
    return $xml
}

function Get-ParagraphAfterQuestion([string]$questionPattern) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -match $questionPattern) {
            return $p.Next()
        }
    }
    throw "Could not locate question paragraph matching: $questionPattern"
}

function Set-SingleAnswerParagraph([string]$questionPattern, [string[]]$runTexts) {
    # Fills the (existing, empty) paragraph right after the question
    # paragraph with the supplied answer text/runs, leaving the blank
    # spacer paragraph that follows untouched.
    $target = Get-ParagraphAfterQuestion $questionPattern
    $xml = New-AnswerParagraphXml $runTexts
    [void]$target.Range.InsertXML($xml)
}

# ------------------------------------------------------------------
# Q8 answer: the two runs ("The mean" + " gestational age ...") are
# merged back into a single run with identical text.
# ------------------------------------------------------------------
$q8Text = "The mean gestational age at birth of babies born to smoking mothers is 38.95, compared with 39.45 for non-smoking mothers. Similarly, the median is 39 for smoking mothers and 40 for non-smoking mothers. This suggests that the pregnancy period is shorter for smoking mothers than for non-smoking mothers."
[void]$d.Content.Find.Execute($q8Text, $true, $false, $false, $false, $false, $true, 1, $false, $q8Text, 2)

# ------------------------------------------------------------------
# Q11 answer: "Yes"
# ------------------------------------------------------------------
Set-SingleAnswerParagraph "Q11\." @("Yes")

# ------------------------------------------------------------------
# Q12 answer: "0.37"
# ------------------------------------------------------------------
Set-SingleAnswerParagraph "Q12\." @("0.37")

# ------------------------------------------------------------------
# Q13 answer: "0.00" — here the question is followed by an empty
# paragraph AND a paragraph that only holds a single space run; the
# first becomes the answer, the second is cleared to a plain empty
# paragraph.
# ------------------------------------------------------------------
$q13Target = Get-ParagraphAfterQuestion "Q13\."
$q13Spacer = $q13Target.Next()
[void]$q13Target.Range.InsertXML((New-AnswerParagraphXml @("0.00")))
[void]$q13Spacer.Range.InsertXML('<w:p ' + $wNs + '></w:p>')

# ------------------------------------------------------------------
# Q14 answer: three paragraphs discussing mode/median/mean and the
# resulting skew, replacing the single empty placeholder paragraph.
# ------------------------------------------------------------------
$q14Target = Get-ParagraphAfterQuestion "Q14\."
$q14Xml = (New-AnswerParagraphXml @("Mode = 2.65 ", "– Me", "di", "an = 3.39 – Mean = 3.51")) +
          (New-AnswerParagraphXml @("Mode &lt; Me", "di", "an &lt; ", "Mean")) +
          (New-AnswerParagraphXml @("They are positively ", "skewed."))
[void]$q14Target.Range.InsertXML($q14Xml)

# ------------------------------------------------------------------
# Q15 answer: "Yes"
# ------------------------------------------------------------------
Set-SingleAnswerParagraph "Q15\." @("Yes")

# ------------------------------------------------------------------
# Q16 answer: "0.94"
# ------------------------------------------------------------------
Set-SingleAnswerParagraph "Q16\." @("0.94")

Write-Output "Answers for Q8 (merge), Q11, Q12, Q13, Q14, Q15, Q16 applied."
